$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 209; this shifts the existing rows 209-313
# down to 210-314 (matching the dimension growing from A1:R313 to A1:R314).
$ws.Rows.Item(209).Insert()

# Populate the freshly inserted row 209 with the new weekly record.
$ws.Cells.Item(209, 1).Value = 6
$ws.Cells.Item(209, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(209, 3).Value = "Metropolitana"
$ws.Cells.Item(209, 4).Value = 45016
$ws.Cells.Item(209, 5).Value = 13
$ws.Cells.Item(209, 6).Value = 100112029
$ws.Cells.Item(209, 7).Value = "Orégano"
$ws.Cells.Item(209, 8).Value = "Sin especificar"
$ws.Cells.Item(209, 9).Value = "Primera"
$ws.Cells.Item(209, 10).Value = 44
$ws.Cells.Item(209, 11).Value = 17000
$ws.Cells.Item(209, 12).Value = 18000
$ws.Cells.Item(209, 13).Value = 17477
$ws.Cells.Item(209, 14).Value = "$/docena de atados"
$ws.Cells.Item(209, 15).Value = "Región Metropolitana"
$ws.Cells.Item(209, 16).Value = 5826
$ws.Cells.Item(209, 17).Value = 3
$ws.Cells.Item(209, 18).Value = "Hortaliza"

# Keep the date column's date-formatted style consistent with the other
# rows in this column (the Insert() already copies formatting down from the
# row above, but set it explicitly to be safe).
$ws.Cells.Item(209, 4).NumberFormat = $ws.Cells.Item(210, 4).NumberFormat
